$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StudentLogin")

# --- Data shift -------------------------------------------------------
# Row 4 (SaviTalent / Innovapath9 / GH, linked to mailto:raj@abc.com with
# display text "raj@abc.com") moves down to become the new row 7; rows
# 5-7 (GP / TW / Normal rows) shift up to become rows 4-6.
$ws.Rows(4).Delete()

# --- Hyperlinks ---------------------------------------------------------
# The simulated engine does not re-anchor hyperlink ranges when rows are
# deleted, so rebuild the hyperlink collection from scratch to match the
# new layout (and re-apply the "Hyperlink" cell style that Hyperlinks.Add
# resets each target cell to, so column A keeps style index 1).
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:savitha.ip9@gmail.com")
$ws.Range("A2").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:savitha.ip9@gmail.com")
$ws.Range("A3").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("A6"), "mailto:raj@abc.com")
$ws.Range("A6").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("A7"), "mailto:raj@abc.com", "", "", "raj@abc.com")
$ws.Range("A7").Style = "Hyperlink"

# Write the real row-7 content after the hyperlink is in place, so the
# cell text becomes "SaviTalent" while the hyperlink keeps its own
# display text of "raj@abc.com" (the two differ, as in the source file).
$ws.Range("A7").Value = "SaviTalent"
$ws.Range("B7").Value = "Innovapath9"
$ws.Range("C7").Value = "GH"

# --- View state ---------------------------------------------------------
# Zoom to 205% and select A8 (a blank cell just below the table), which
# also drops the previous topLeftCell scroll anchor.
$excel.ActiveWindow.Zoom = 205
$ws.Range("A8").Select() | Out-Null

# --- Page setup -----------------------------------------------------------
# Touch PageSetup so a <pageSetup orientation="portrait"/> element is
# written out for the sheet.
$ws.PageSetup.Orientation = 1
